$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 228, shifting existing rows 228-232 down to 229-233
$ws.Rows.Item(228).Insert()

# Fill in the new row 228 with the new weekly price entry.
# Columns A,B,C,E,F,G,H,O,R are constant for this market/product, copy them from the row below (now 229)
$ws.Range("A228").Value2 = $ws.Range("A229").Value2
$ws.Range("B228").Value2 = $ws.Range("B229").Value2
$ws.Range("C228").Value2 = $ws.Range("C229").Value2
$ws.Range("E228").Value2 = $ws.Range("E229").Value2
$ws.Range("F228").Value2 = $ws.Range("F229").Value2
$ws.Range("G228").Value2 = $ws.Range("G229").Value2
$ws.Range("H228").Value2 = $ws.Range("H229").Value2
$ws.Range("O228").Value2 = $ws.Range("O229").Value2
$ws.Range("R228").Value2 = $ws.Range("R229").Value2

# Copy the date cell style from row 229's D cell so the new D228 keeps the date format
$ws.Range("D229").Copy($ws.Range("D228"))

$ws.Range("D228").Value2 = 45239
$ws.Range("I228").Value2 = "Primera"
$ws.Range("J228").Value2 = 500
$ws.Range("K228").Value2 = 12000
$ws.Range("L228").Value2 = 13000
$ws.Range("M228").Value2 = 12500
$ws.Range("N228").Value2 = "`$/caja 70 unidades"
$ws.Range("P228").Value2 = 179
$ws.Range("Q228").Value2 = 70
